$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.039.85"
$ws.Range("E2").Value = "  -4.41%  "
$ws.Range("D3").Value = "3.281.42"
$ws.Range("E3").Value = "  -4.98%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'558.60"
$ws.Range("E5").Value = "  -2.97%  "
$ws.Range("D6").Value = "'185.35"
$ws.Range("E6").Value = "  -2.58%  "
$ws.Range("E8").Value = "  -2.36%  "
$ws.Range("D9").Value = "3.275.49"
$ws.Range("E9").Value = "  -4.92%  "
$ws.Range("E10").Value = "  -7.41%  "
$ws.Range("D11").Value = "'0.588"
$ws.Range("E11").Value = "  -4.44%  "
$ws.Range("D12").Value = "'47.48"
$ws.Range("E12").Value = "  -7.21%  "
$ws.Range("E13").Value = "  -5.76%  "
$ws.Range("E14").Value = "  -5.14%  "
$ws.Range("D15").Value = "'629.21"
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("D16").Value = "3.809.98"
$ws.Range("E16").Value = "  -4.63%  "
$ws.Range("D17").Value = "66.027.86"
$ws.Range("E17").Value = "  -4.23%  "
$ws.Range("D18").Value = "'17.88"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("E19").Value = "  -3.22%  "
$ws.Range("D20").Value = "3.279.88"
$ws.Range("E20").Value = "  -4.56%  "
$ws.Range("D21").Value = "'11.39"
$ws.Range("E21").Value = "  -6.66%  "
$ws.Range("E22").Value = "  -3.27%  "
$ws.Range("D23").Value = "'17.98"
$ws.Range("E23").Value = "  +1.52%  "
$ws.Range("D24").Value = "'106.42"
$ws.Range("E24").Value = "  +8.35%  "
$ws.Range("D25").Value = "'4.94"
$ws.Range("E25").Value = "  -6.74%  "
$ws.Range("D26").Value = "'3.97"
$ws.Range("E26").Value = "  -6.69%  "
$ws.Range("E27").Value = "  -5.72%  "
$ws.Range("D28").Value = "'9.65"
$ws.Range("E28").Value = "  -1.94%  "
$ws.Range("D29").Value = "'8.73"
$ws.Range("E29").Value = "  -5.19%  "
$ws.Range("E30").Value = "  -5.19%  "
$ws.Range("D31").Value = "'4.04"
$ws.Range("E31").Value = "  -6.04%  "
$ws.Range("D32").Value = "'6.28"
$ws.Range("E32").Value = "  -5.82%  "
$ws.Range("D33").Value = "'11.05"
$ws.Range("E33").Value = "  -4.10%  "
$ws.Range("E34").Value = "  -3.00%  "
$ws.Range("D35").Value = "'539.50"
$ws.Range("E35").Value = "  +8.44%  "
$ws.Range("D36").Value = "'57.48"
$ws.Range("E36").Value = "  -5.71%  "
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("D38").Value = "3.690.28"
$ws.Range("E38").Value = "  +0.90%  "
$ws.Range("D39").Value = "'3.44"
$ws.Range("E40").Value = "  -6.62%  "
$ws.Range("D41").Value = "'0.132"
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").Value = "'2.76"
$ws.Range("E42").Value = "  -5.38%  "
$ws.Range("D43").Value = "'3.36"
$ws.Range("E43").Value = "  -3.97%  "
$ws.Range("D44").Value = "'32.69"
$ws.Range("E44").Value = "  -4.12%  "
$ws.Range("D45").Value = "'0.339"
$ws.Range("E45").Value = "  -7.96%  "
$ws.Range("E46").Value = "  -1.99%  "
$ws.Range("D47").Value = "'0.0415"
$ws.Range("E47").Value = "  -4.59%  "
$ws.Range("E48").Value = "  -6.28%  "
$ws.Range("E49").Value = "  -3.15%  "
$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("E51").Value = "  +2.24%  "

Write-Host "Applied 81 cell updates"
